# Auto-generated Excel COM-interop edit script
# Applies updated profit-calculation values to the Coeurl_Profits workbook
# (scheduled-runner price refresh across all DoH/DoL job sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1662.1666
$ws.Range("I31").Value = 1662.1666
$ws.Range("K31").Value = 4986.4998
$ws.Range("M31").Value = -4756.4998
$ws.Range("H39").Value = 75.85714
$ws.Range("I39").Value = 66.40000000000001
$ws.Range("J39").Value = 99.5
$ws.Range("K39").Value = 199.2
$ws.Range("L39").Value = 298.5
$ws.Range("M39").Value = 96.79999999999998
$ws.Range("N39").Value = -890.5
$ws.Range("H64").Value = 7025.778
$ws.Range("I64").Value = 4106.1
$ws.Range("J64").Value = 10675.375
$ws.Range("K64").Value = 4106.1
$ws.Range("L64").Value = 10675.375
$ws.Range("M64").Value = -3858.1
$ws.Range("N64").Value = -11171.375
$ws.Range("H67").Value = 7025.778
$ws.Range("I67").Value = 4106.1
$ws.Range("J67").Value = 10675.375
$ws.Range("K67").Value = 4106.1
$ws.Range("L67").Value = 10675.375
$ws.Range("M67").Value = -3248.1
$ws.Range("N67").Value = -12391.375
$ws.Range("H69").Value = 8588.799999999999
$ws.Range("J69").Value = 8588.799999999999
$ws.Range("L69").Value = 25766.4
$ws.Range("N69").Value = -27514.4
$ws.Range("H70").Value = 4017.842
$ws.Range("I70").Value = 4819.9165
$ws.Range("J70").Value = 2642.8572
$ws.Range("K70").Value = 14459.7495
$ws.Range("L70").Value = 7928.571599999999
$ws.Range("M70").Value = -14189.7495
$ws.Range("N70").Value = -8468.571599999999
$ws.Range("H72").Value = 8588.799999999999
$ws.Range("J72").Value = 8588.799999999999
$ws.Range("L72").Value = 77299.2
$ws.Range("N72").Value = -86035.2
$ws.Range("H73").Value = 4017.842
$ws.Range("I73").Value = 4819.9165
$ws.Range("J73").Value = 2642.8572
$ws.Range("K73").Value = 14459.7495
$ws.Range("L73").Value = 7928.571599999999
$ws.Range("M73").Value = -13523.7495
$ws.Range("N73").Value = -9800.571599999999
$ws.Range("H74").Value = 7385.619
$ws.Range("I74").Value = 5283.222
$ws.Range("J74").Value = 20000
$ws.Range("K74").Value = 5283.222
$ws.Range("L74").Value = 20000
$ws.Range("M74").Value = -4347.222
$ws.Range("N74").Value = -21872
$ws.Range("H77").Value = 7385.619
$ws.Range("I77").Value = 5283.222
$ws.Range("J77").Value = 20000
$ws.Range("K77").Value = 26416.11
$ws.Range("L77").Value = 100000
$ws.Range("M77").Value = -21736.11
$ws.Range("N77").Value = -109360
$ws.Range("H94").Value = 608.1818
$ws.Range("I94").Value = 608.1818
$ws.Range("K94").Value = 608.1818
$ws.Range("M94").Value = -157.1818
$ws.Range("H98").Value = 1483.5
$ws.Range("I98").Value = 1203.1666
$ws.Range("K98").Value = 1203.1666
$ws.Range("M98").Value = 294.8334
$ws.Range("H112").Value = 45361.64
$ws.Range("I112").Value = 2075
$ws.Range("J112").Value = 53606.715
$ws.Range("K112").Value = 6225
$ws.Range("L112").Value = 160820.145
$ws.Range("M112").Value = -5117
$ws.Range("N112").Value = -163036.145
$ws.Range("H118").Value = 2553.8462
$ws.Range("I118").Value = 458.2857
$ws.Range("K118").Value = 1374.8571
$ws.Range("M118").Value = 282.1428999999998
$ws.Range("H122").Value = 1483.5
$ws.Range("I122").Value = 1203.1666
$ws.Range("K122").Value = 3609.4998
$ws.Range("M122").Value = -1159.4998
$ws.Range("H135").Value = 1064.6061
$ws.Range("I135").Value = 875.2258
$ws.Range("K135").Value = 7877.032200000001
$ws.Range("M135").Value = -5342.032200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1624
$ws.Range("I2").Value = 251.66667
$ws.Range("K2").Value = 251.66667
$ws.Range("M2").Value = -138.66667
$ws.Range("H33").Value = 26500
$ws.Range("I33").Value = 32000
$ws.Range("J33").Value = 15500
$ws.Range("K33").Value = 32000
$ws.Range("L33").Value = 15500
$ws.Range("M33").Value = -31671
$ws.Range("N33").Value = -16158
$ws.Range("H39").Value = 898.75
$ws.Range("I39").Value = 898.75
$ws.Range("K39").Value = 898.75
$ws.Range("M39").Value = -378.75
$ws.Range("H97").Value = 1303.0541
$ws.Range("I97").Value = 970.1212
$ws.Range("K97").Value = 970.1212
$ws.Range("M97").Value = -474.1212
$ws.Range("H116").Value = 1624
$ws.Range("I116").Value = 251.66667
$ws.Range("K116").Value = 251.66667
$ws.Range("M116").Value = 2042.33333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1624
$ws.Range("I3").Value = 251.66667
$ws.Range("K3").Value = 251.66667
$ws.Range("M3").Value = -137.66667
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 64552.25
$ws.Range("I31").Value = 78605.16
$ws.Range("J31").Value = 3656.3333
$ws.Range("K31").Value = 78605.16
$ws.Range("L31").Value = 3656.3333
$ws.Range("M31").Value = -78310.16
$ws.Range("N31").Value = -4246.3333
$ws.Range("H34").Value = 64552.25
$ws.Range("I34").Value = 78605.16
$ws.Range("J34").Value = 3656.3333
$ws.Range("K34").Value = 78605.16
$ws.Range("L34").Value = 3656.3333
$ws.Range("M34").Value = -78403.16
$ws.Range("N34").Value = -4060.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 24991.5
$ws.Range("I56").Value = 24991.5
$ws.Range("K56").Value = 24991.5
$ws.Range("M56").Value = -24461.5
$ws.Range("H68").Value = 12500834
$ws.Range("J68").Value = 16667579
$ws.Range("L68").Value = 50002737
$ws.Range("N68").Value = -50004359
$ws.Range("H71").Value = 12500834
$ws.Range("J71").Value = 16667579
$ws.Range("L71").Value = 150008211
$ws.Range("N71").Value = -150016323
$ws.Range("H121").Value = 776
$ws.Range("I121").Value = 290
$ws.Range("J121").Value = 1359.2
$ws.Range("K121").Value = 870
$ws.Range("L121").Value = 4077.6
$ws.Range("M121").Value = 440
$ws.Range("N121").Value = -6697.6
$ws.Range("H122").Value = 1207.4286
$ws.Range("I122").Value = 928.6429000000001
$ws.Range("J122").Value = 1765
$ws.Range("K122").Value = 8357.786100000001
$ws.Range("L122").Value = 15885
$ws.Range("M122").Value = -5907.786100000001
$ws.Range("N122").Value = -20785
$ws.Range("H132").Value = 1402.3438
$ws.Range("I132").Value = 1306.7307
$ws.Range("J132").Value = 1816.6666
$ws.Range("K132").Value = 11760.5763
$ws.Range("L132").Value = 16349.9994
$ws.Range("M132").Value = -9230.576300000001
$ws.Range("N132").Value = -21409.9994
$ws.Range("H140").Value = 1574.4
$ws.Range("I140").Value = 1574.4
$ws.Range("K140").Value = 4723.200000000001
$ws.Range("M140").Value = 456.7999999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2721.2222
$ws.Range("I122").Value = 2699
$ws.Range("J122").Value = 2899
$ws.Range("K122").Value = 8097
$ws.Range("L122").Value = 8697
$ws.Range("M122").Value = -5647
$ws.Range("N122").Value = -13597

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4124.1665
$ws.Range("I40").Value = 3191.1538
$ws.Range("K40").Value = 3191.1538
$ws.Range("M40").Value = -3055.1538
$ws.Range("H47").Value = 50000
$ws.Range("J47").Value = 50000
$ws.Range("L47").Value = 50000
$ws.Range("N47").Value = -50980
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50466
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H82").Value = 6663.0835
$ws.Range("I82").Value = 9811.385
$ws.Range("K82").Value = 9811.385
$ws.Range("M82").Value = -9450.385
$ws.Range("H85").Value = 6663.0835
$ws.Range("I85").Value = 9811.385
$ws.Range("K85").Value = 9811.385
$ws.Range("M85").Value = -8563.385
$ws.Range("H100").Value = 335266.66
$ws.Range("J100").Value = 502024.5
$ws.Range("L100").Value = 502024.5
$ws.Range("N100").Value = -503106.5
$ws.Range("H132").Value = 31506.514
$ws.Range("I132").Value = 36887
$ws.Range("K132").Value = 110661
$ws.Range("M132").Value = -108131

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 29499.5
$ws.Range("I9").Value = 44500.5
$ws.Range("J9").Value = 14498.5
$ws.Range("K9").Value = 44500.5
$ws.Range("L9").Value = 14498.5
$ws.Range("M9").Value = -44360.5
$ws.Range("N9").Value = -14778.5
$ws.Range("H62").Value = 6900
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -13248
$ws.Range("H65").Value = 6900
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -66240
$ws.Range("H81").Value = 7881.727
$ws.Range("I81").Value = 12744.223
$ws.Range("J81").Value = 4515.385
$ws.Range("K81").Value = 25488.446
$ws.Range("L81").Value = 9030.77
$ws.Range("M81").Value = -24427.446
$ws.Range("N81").Value = -11152.77
$ws.Range("H84").Value = 7881.727
$ws.Range("I84").Value = 12744.223
$ws.Range("J84").Value = 4515.385
$ws.Range("K84").Value = 127442.23
$ws.Range("L84").Value = 45153.85000000001
$ws.Range("M84").Value = -122138.23
$ws.Range("N84").Value = -55761.85000000001
